$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4741
$ws1.Range("F3").Value = 2742
$ws1.Range("F5").Value = 2780
$ws1.Range("F7").Value = 1952
$ws1.Range("F12").Value = 233
$ws1.Range("F13").Value = 407
$ws1.Range("F14").Value = 1062
$ws1.Range("F15").Value = 298
$ws1.Range("F18").Value = 531
$ws1.Range("F19").Value = 531
$ws1.Range("F22").Value = 657
$ws1.Range("F26").Value = 511
$ws1.Range("F27").Value = 15
$ws1.Range("F29").Value = 1518
$ws1.Range("F30").Value = 331
$ws1.Range("F32").Value = 1472
$ws1.Range("F33").Value = 154
$ws1.Range("F34").Value = 2327
$ws1.Range("F35").Value = 386
$ws1.Range("F36").Value = 25
$ws1.Range("F37").Value = 606
$ws1.Range("F39").Value = 59
$ws1.Range("F41").Value = 780
$ws1.Range("F42").Value = 1481
$ws1.Range("F43").Value = 207
$ws1.Range("F45").Value = 492
$ws1.Range("F46").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4741
$ws4.Range("F3").Value = 2742
$ws4.Range("F4").Value = 2780
$ws4.Range("F10").Value = 233
$ws4.Range("F11").Value = 407
$ws4.Range("F12").Value = 1062
$ws4.Range("F13").Value = 298
$ws4.Range("F16").Value = 531
$ws4.Range("F17").Value = 531
$ws4.Range("F19").Value = 657
$ws4.Range("F26").Value = 511
$ws4.Range("F28").Value = 1518
$ws4.Range("F29").Value = 331
$ws4.Range("F33").Value = 2327
$ws4.Range("F34").Value = 386
$ws4.Range("F37").Value = 25
$ws4.Range("F39").Value = 606
$ws4.Range("F41").Value = 59
$ws4.Range("F43").Value = 780
$ws4.Range("F44").Value = 1481
$ws4.Range("F46").Value = 207
$ws4.Range("F47").Value = 492
